$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Lowell"
$ws.Range("B1").Value = "5ouolkqhu0x83vp"
$ws.Range("A2").Value = "Otto"
$ws.Range("B2").Value = "1w9uld5iyb0k0"
$ws.Range("A3").Value = "Lorenza"
$ws.Range("B3").Value = "2he96xy96b"
